# Day7 - completed all the major parts, left with column subtract operation
# and group by labels on the bottom of the page.
#
# Adds three "recal" columns (AG:AI -> budget_recal / april_reforecast_recal /
# actuals_recal) mirroring the Budget/April-Reforecast/Actuals FY totals,
# fills in previously-blank grouping columns (B:F, G:J, S:T) that were left
# as empty placeholders by the prior export, and appends two new summary
# rows ("recal" and "variance") at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header rows: new recal columns
# ---------------------------------------------------------------------
$ws.Range("AG1").Value = "budget_recal"
$ws.Range("AH1").Value = "april_reforecast_recal"
$ws.Range("AI1").Value = "actuals_recal"

$ws.Range("AG2").Value = "budget_recal"
$ws.Range("AH2").Value = "april_reforecast_recal"
$ws.Range("AI2").Value = "actuals_recal"

# ---------------------------------------------------------------------
# Row 3 (Study 1 / Cost Type A)
# ---------------------------------------------------------------------
$ws.Range("G3:J3").Value = 0
$ws.Range("S3:T3").Value = 0
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 150.4579185173292
$ws.Range("AI3").Value = 157.0627923229017

# ---------------------------------------------------------------------
# Row 4 (Study 1 / Cost Type B)
# ---------------------------------------------------------------------
$ws.Range("B4").Value = "Development"
$ws.Range("C4").Value = "Clinical"
$ws.Range("E4").Value = "Study 1"
$ws.Range("G4:J4").Value = 0
$ws.Range("S4:T4").Value = 0
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 826.1550247118
$ws.Range("AI4").Value = 6865.47324

# ---------------------------------------------------------------------
# Row 5 (54321 Total)
# ---------------------------------------------------------------------
$ws.Range("B5").Value = "Development"
$ws.Range("C5").Value = "Clinical"
$ws.Range("E5").Value = "Study 1"
$ws.Range("F5").Value = "B"
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 976.6129432291293
$ws.Range("AI5").Value = 7022.536032322902

# ---------------------------------------------------------------------
# Row 6 (Study 2 / Cost Type A)
# ---------------------------------------------------------------------
$ws.Range("B6").Value = "Development"
$ws.Range("C6").Value = "Clinical"
$ws.Range("G6:J6").Value = 0
$ws.Range("S6:T6").Value = 0
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 291.2865305604031
$ws.Range("AI6").Value = 180.0596631326561

# ---------------------------------------------------------------------
# Row 7 (Study 2 / Cost Type B)
# ---------------------------------------------------------------------
$ws.Range("B7").Value = "Development"
$ws.Range("C7").Value = "Clinical"
$ws.Range("D7").Value = 65432
$ws.Range("E7").Value = "Study 2"
$ws.Range("G7:J7").Value = 0
$ws.Range("S7:T7").Value = 0
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 1912.072160284128
$ws.Range("AI7").Value = 158.54728

# ---------------------------------------------------------------------
# Row 8 (65432 Total)
# ---------------------------------------------------------------------
$ws.Range("B8").Value = "Development"
$ws.Range("C8").Value = "Clinical"
$ws.Range("E8").Value = "Study 2"
$ws.Range("F8").Value = "B"
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 2203.358690844531
$ws.Range("AI8").Value = 338.6069431326561

# ---------------------------------------------------------------------
# Row 9 (Study 3 / Cost Type A)
# ---------------------------------------------------------------------
$ws.Range("B9").Value = "Development"
$ws.Range("C9").Value = "Clinical"
$ws.Range("G9:J9").Value = 0
$ws.Range("S9:T9").Value = 0
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 1354.208794338037
$ws.Range("AI9").Value = 1085.326284925343

# ---------------------------------------------------------------------
# Row 10 (87653 Total)
# ---------------------------------------------------------------------
$ws.Range("B10").Value = "Development"
$ws.Range("C10").Value = "Clinical"
$ws.Range("E10").Value = "Study 3"
$ws.Range("F10").Value = "A"
$ws.Range("AG10").Value = 0
$ws.Range("AH10").Value = 1354.208794338037
$ws.Range("AI10").Value = 1085.326284925343

# ---------------------------------------------------------------------
# Row 11 (Clinical Total)
# ---------------------------------------------------------------------
$ws.Range("B11").Value = "Development"
$ws.Range("D11").Value = "87653 Total"
$ws.Range("E11").Value = "Study 3"
$ws.Range("F11").Value = "A"
$ws.Range("AG11").Value = 0
$ws.Range("AH11").Value = 4534.180428411697
$ws.Range("AI11").Value = 8446.469260380902

# ---------------------------------------------------------------------
# Row 12 (Development Total)
# ---------------------------------------------------------------------
$ws.Range("C12").Value = "Clinical Total"
$ws.Range("D12").Value = "87653 Total"
$ws.Range("E12").Value = "Study 3"
$ws.Range("F12").Value = "A"
$ws.Range("AG12").Value = 0
$ws.Range("AH12").Value = 4534.180428411697
$ws.Range("AI12").Value = 8446.469260380902

# ---------------------------------------------------------------------
# Row 13 (Grand Total)
# ---------------------------------------------------------------------
$ws.Range("C13").Value = "Clinical Total"
$ws.Range("D13").Value = "87653 Total"
$ws.Range("E13").Value = "Study 3"
$ws.Range("F13").Value = "A"
$ws.Range("AG13").Value = 0
$ws.Range("AH13").Value = 4534.180428411697
$ws.Range("AI13").Value = 8446.469260380902

# ---------------------------------------------------------------------
# Row 14: new "recal" summary row (same figures as the Grand Total row,
# recalculated) - copy formatting of the label cell from A13 so the new
# label keeps the bold/centered/bordered look used throughout column A.
# ---------------------------------------------------------------------
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = "recal"

$ws.Range("G14:K14").Value = 0
$ws.Range("L14").Value = 1417.598704893892
$ws.Range("M14").Value = 604.61526963542
$ws.Range("N14").Value = 2226.315374733448
$ws.Range("O14").Value = 285.6510791489375
$ws.Range("P14").Value = 4534.180428411697
$ws.Range("Q14").Value = 999.1874956821435
$ws.Range("R14").Value = 7447.281764698757
$ws.Range("S14:T14").Value = 0
$ws.Range("U14").Value = 8446.469260380902
$ws.Range("V14").Value = 999.1874956821435
$ws.Range("W14").Value = 7447.281764698757
$ws.Range("X14:Y14").Value = 0
$ws.Range("Z14").Value = 8446.469260380902
$ws.Range("AA14").Value = -418.4112092117485
$ws.Range("AB14").Value = 6842.666495063337
$ws.Range("AC14").Value = -2226.315374733448
$ws.Range("AD14").Value = -285.6510791489375
$ws.Range("AE14").Value = 3912.288831969205

# ---------------------------------------------------------------------
# Row 15: new "variance" summary row (placeholder row for the upcoming
# column-subtract operation mentioned in the commit message - all zero
# for now).
# ---------------------------------------------------------------------
$ws.Range("A13").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = "variance"

$ws.Range("G15:AE15").Value = 0

Write-Host "edit complete"
